$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text, matching original inline-string cell type
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.233.50"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "1.912.83"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "0.7398"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").Value = "243.25"
$ws.Range("E6").Value = "  -2.24%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "0.3144"
$ws.Range("E8").Value = "  -1.86%  "

$ws.Range("D9").Value = "27.12"
$ws.Range("E9").Value = "  -3.57%  "

$ws.Range("D10").Value = "0.06968"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("D11").Value = "0.7803"
$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").Value = "0.07974"
$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("D13").Value = "1.926.73"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").Value = "5.276"
$ws.Range("E14").Value = "  -1.93%  "

$ws.Range("D15").Value = "91.52"
$ws.Range("E15").Value = "  -3.14%  "

$ws.Range("D16").Value = "30.317.54"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "14.28"
$ws.Range("E17").Value = "  -2.29%  "

$ws.Range("D18").Value = "245.40"
$ws.Range("E18").Value = "  -2.87%  "

$ws.Range("D19").Value = "5.820"
$ws.Range("E19").Value = "  +0.60%  "

$ws.Range("D20").Value = "0.000007818"
$ws.Range("E20").Value = "  -2.73%  "

$ws.Range("D21").Value = "2.208.05"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "6.629"
$ws.Range("E24").Value = "  -2.86%  "

$ws.Range("D25").Value = "9.405"
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").Value = "164.87"
$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").Value = "19.00"
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("D28").Value = "0.1274"
$ws.Range("E28").Value = "  -3.14%  "

$ws.Range("D29").Value = "2.132"
$ws.Range("E29").Value = "  -8.53%  "

$ws.Range("D30").Value = "1.352"
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").Value = "1.545"
$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").Value = "4.324"
$ws.Range("E32").Value = "  -2.50%  "

$ws.Range("D33").Value = "4.085"
$ws.Range("E33").Value = "  -1.56%  "

$ws.Range("D34").Value = "0.05198"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").Value = "1.298"
$ws.Range("E35").Value = "  +1.11%  "

$ws.Range("D36").Value = "0.7526"
$ws.Range("E36").Value = "  +0.46%  "

$ws.Range("D37").Value = "2.759"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("D38").Value = "0.01941"
$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("D39").Value = "2.793"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "6.394"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").Value = "76.10"
$ws.Range("E41").Value = "  -2.52%  "

$ws.Range("D42").Value = "0.4492"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("D43").Value = "1.948"
$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").Value = "0.8336"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").Value = "7.673"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D47").Value = "101.38"
$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("D48").Value = "9.909"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").Value = "2.125.32"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").Value = "37.04"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").Value = "0.1219"
$ws.Range("E51").Value = "  +2.32%  "
